# Fill in the three description/grade pairs (Math, English, Jewish-History
# tables) that were left empty in the source certificate template. The
# first table (Torah) is intentionally left untouched, matching the diff.

$d = $word.ActiveDocument

function Get-ParagraphAtStart($doc, $targetStart) {
    $match = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Start -eq $targetStart) {
            $match = $p
        }
    }
    return $match
}

function Find-NextRange($doc, $searchText, $fromPos) {
    $docEnd = $doc.Content.End
    $r = $doc.Range($fromPos, $docEnd)
    $found = $r.Find.Execute($searchText)
    if (-not $found) {
        return $null
    }
    return $r
}

function Set-SubjectDescriptionAndGrade($doc, $subjectText, $descriptionText, $gradeText) {
    # Locate the subject-name paragraph (e.g. the "Math" heading); the very
    # next paragraph in document order is the (currently empty) free-text
    # description cell belonging to the merged right-hand column.
    $subjRange = $doc.Range(0, $doc.Content.End)
    $found = $subjRange.Find.Execute($subjectText)
    if (-not $found) {
        throw "subject not found: $subjectText"
    }
    $subjPara = Get-ParagraphAtStart $doc $subjRange.Start
    if ($subjPara -eq $null) {
        throw "subject paragraph not found: $subjectText"
    }
    $descPara = $subjPara.Next()
    $descPara.Range.Text = $descriptionText

    # Locate the nearest "ציון:" label that follows the subject; the
    # paragraph right after it is the (currently empty) grade cell.
    $labelRange = Find-NextRange $doc "ציון:" $subjRange.End
    if ($labelRange -eq $null) {
        throw "grade label not found for: $subjectText"
    }
    $labelPara = Get-ParagraphAtStart $doc $labelRange.Start
    if ($labelPara -eq $null) {
        throw "grade label paragraph not found for: $subjectText"
    }
    $gradePara = $labelPara.Next()
    $gradePara.Range.Text = $gradeText
}

Set-SubjectDescriptionAndGrade $d "מתמטיקה" "במחצית למדנו משוואות ב2 נעלמים, פיתחנו כמה שיטות לבעיה זו,בנוסף התעסקנו בבעיות תנועה וזמן ולמדנו איך להתמודד מול זאת`nהיה לנו הספקים מעולים!`nשפרה את ילדה מקסימה, שיהיה לך הרבה הצלחה בהמשך! " "89"
Set-SubjectDescriptionAndGrade $d "אנגלית" "במחצית זאת התמקדנו על הבנה חזקה של הטקסטים ולמדנו את השיטות להבנת הנקרא, חזרנו על שאלות חוזרות ופיתחנו שיטות קלות לפיתרתן.`nשפרה את ילדה נהדרת, הרבה הצלחה!" "96"
Set-SubjectDescriptionAndGrade $d "תולדות ישראל" "במחצית זאת למדנו על גדולי ישראל בכל מיני יבשות, על המצב של היהודים בתקופות שלטון שונות,`nשפרה הרבה הצלחה!" "90"
